$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. All D/E column cells in this sheet are
# stored as literal text (prices / percentages with formatting such as
# trailing zeros or a "%" suffix) - force text number-format before writing
# so Excel does not auto-coerce the string into a numeric/percentage value,
# then restore the default "Normal" style so no stray formatting is left
# behind on cells that originally had none.
$updates = @{
    "D2" = "290.88"
    "E2" = "-3.34%"
    "D3" = "30.65"
    "E3" = "-6.32%"
    "D4" = "4.944"
    "E4" = "0.13%"
    "D5" = "0.07205"
    "E5" = "-6.56%"
    "D6" = "1.807"
    "E6" = "-8.13%"
    "D7" = "7.678"
    "E7" = "-2.00%"
    "D8" = "3.768"
    "E8" = "-0.82%"
    "D9" = "0.8966"
    "E9" = "-2.55%"
    "D10" = "0.1653"
    "E10" = "-5.71%"
    "D11" = "0.07736"
    "E11" = "-0.60%"
    "D12" = "0.08068"
    "E12" = "-6.09%"
    "D13" = "0.03059"
    "E13" = "-3.90%"
    "D14" = "0.1002"
    "E14" = "-0.10%"
    "D15" = "0.001495"
    "E15" = "-1.85%"
    "D16" = "0.005731"
    "E16" = "-1.49%"
    "D18" = "3.465"
    "E18" = "0.18%"
    "E19" = "-3.37%"
    "E20" = "-0.97%"
    "D21" = "0.1299"
    "E21" = "-2.07%"
    "D22" = "4.041"
    "E22" = "-5.37%"
    "D23" = "0.2391"
    "E23" = "20.03%"
    "D24" = "0.04496"
    "E24" = "-0.60%"
    "D25" = "0.001217"
    "E25" = "-0.42%"
    "D26" = "0.004002"
    "E26" = "-9.39%"
    "D27" = "0.0001252"
    "E27" = "0.00%"
    "D39" = "0.01583"
    "E39" = "-6.95%"
    "D40" = "0.04407"
    "E40" = "-6.08%"
    "D41" = "0.007247"
    "E41" = "-3.21%"
    "D42" = "0.009914"
    "D43" = "0.1307"
    "E43" = "-3.40%"
    "D44" = "0.002009"
    "E44" = "-13.90%"
    "D45" = "0.009506"
    "E45" = "-9.71%"
    "D46" = "0.00005948"
    "E46" = "-4.64%"
    "D47" = "0.00000000751"
    "E47" = "0.02%"
    "E48" = "173.65%"
    "D49" = "0.003004"
    "E49" = "-3.24%"
    "D50" = "0.00002103"
    "E50" = "0.02%"
    "D51" = "0.0002003"
    "E51" = "0.02%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Write-Output "Updated $($updates.Count) cells"
